$d = $word.ActiveDocument

$d.Content.Find.Execute("10 seconds", $false, $false, $false, $false, $false, $true, 1, $false, "7 seconds", 2)
